$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 458
$ws.Range("J33").Value = 599
$ws.Range("L33").Value = 599
$ws.Range("N33").Value = -1057
$ws.Range("H40").Value = 100001700
$ws.Range("I40").Value = 1495
$ws.Range("J40").Value = 125001740
$ws.Range("K40").Value = 1495
$ws.Range("L40").Value = 125001740
$ws.Range("M40").Value = -1320
$ws.Range("N40").Value = -125002090
$ws.Range("H69").Value = 9635.909
$ws.Range("I69").Value = 9166.666999999999
$ws.Range("J69").Value = 9811.875
$ws.Range("K69").Value = 27500.001
$ws.Range("L69").Value = 29435.625
$ws.Range("M69").Value = -26626.001
$ws.Range("N69").Value = -31183.625
$ws.Range("H72").Value = 9635.909
$ws.Range("I72").Value = 9166.666999999999
$ws.Range("J72").Value = 9811.875
$ws.Range("K72").Value = 82500.003
$ws.Range("L72").Value = 88306.875
$ws.Range("M72").Value = -78132.003
$ws.Range("N72").Value = -97042.875
$ws.Range("H113").Value = 3037.818
$ws.Range("I113").Value = 2900.6667
$ws.Range("K113").Value = 2900.6667
$ws.Range("M113").Value = 353.3332999999998
$ws.Range("H133").Value = 116571.14
$ws.Range("J133").Value = 116571.14
$ws.Range("L133").Value = 116571.14
$ws.Range("N133").Value = -126691.14

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 591.25
$ws.Range("I2").Value = 497.3793
$ws.Range("J2").Value = 980.1429000000001
$ws.Range("K2").Value = 497.3793
$ws.Range("L2").Value = 980.1429000000001
$ws.Range("M2").Value = -384.3793
$ws.Range("N2").Value = -1206.1429
$ws.Range("H32").Value = 8707.154
$ws.Range("I32").Value = 8686.953
$ws.Range("K32").Value = 8686.953
$ws.Range("M32").Value = -8399.953
$ws.Range("H74").Value = 1748.3334
$ws.Range("I74").Value = 1621.9524
$ws.Range("K74").Value = 1621.9524
$ws.Range("M74").Value = -747.9523999999999
$ws.Range("H77").Value = 1748.3334
$ws.Range("I77").Value = 1621.9524
$ws.Range("K77").Value = 8109.762
$ws.Range("M77").Value = -3741.762
$ws.Range("H116").Value = 591.25
$ws.Range("I116").Value = 497.3793
$ws.Range("J116").Value = 980.1429000000001
$ws.Range("K116").Value = 497.3793
$ws.Range("L116").Value = 980.1429000000001
$ws.Range("M116").Value = 1796.6207
$ws.Range("N116").Value = -5568.1429
$ws.Range("H132").Value = 1697223.5
$ws.Range("I132").Value = 2201.204
$ws.Range("J132").Value = 10002833
$ws.Range("K132").Value = 6603.612000000001
$ws.Range("L132").Value = 30008499
$ws.Range("M132").Value = -4073.612000000001
$ws.Range("N132").Value = -30013559
$ws.Range("H134").Value = 50749.5
$ws.Range("J134").Value = 50749.5
$ws.Range("L134").Value = 50749.5
$ws.Range("N134").Value = -60889.5
$ws.Range("H140").Value = 29999
$ws.Range("J140").Value = 29999
$ws.Range("L140").Value = 29999
$ws.Range("N140").Value = -40359

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 591.25
$ws.Range("I3").Value = 497.3793
$ws.Range("J3").Value = 980.1429000000001
$ws.Range("K3").Value = 497.3793
$ws.Range("L3").Value = 980.1429000000001
$ws.Range("M3").Value = -383.3793
$ws.Range("N3").Value = -1208.1429
$ws.Range("H22").Value = 1292.3636
$ws.Range("I22").Value = 1324
$ws.Range("K22").Value = 1324
$ws.Range("M22").Value = -1151
$ws.Range("H99").Value = 1925.3334
$ws.Range("I99").Value = 1114.8572
$ws.Range("K99").Value = 1114.8572
$ws.Range("M99").Value = 383.1428000000001
$ws.Range("H134").Value = 11111997
$ws.Range("I134").Value = 997.125
$ws.Range("J134").Value = 100000000
$ws.Range("K134").Value = 2991.375
$ws.Range("L134").Value = 300000000
$ws.Range("M134").Value = -456.375
$ws.Range("N134").Value = -300005070
$ws.Range("H140").Value = 177446.17
$ws.Range("J140").Value = 177446.17
$ws.Range("L140").Value = 177446.17
$ws.Range("N140").Value = -187806.17

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 647.6923
$ws.Range("I22").Value = 392.75
$ws.Range("J22").Value = 1055.6
$ws.Range("K22").Value = 392.75
$ws.Range("L22").Value = 1055.6
$ws.Range("M22").Value = -42.75
$ws.Range("N22").Value = -1755.6
$ws.Range("H51").Value = 25948.334
$ws.Range("J51").Value = 32500
$ws.Range("L51").Value = 32500
$ws.Range("N51").Value = -33972
$ws.Range("H58").Value = 2471
$ws.Range("I58").Value = 1760.8823
$ws.Range("K58").Value = 1760.8823
$ws.Range("M58").Value = -1557.8823
$ws.Range("H61").Value = 25948.334
$ws.Range("J61").Value = 32500
$ws.Range("L61").Value = 32500
$ws.Range("N61").Value = -33196
$ws.Range("H62").Value = 19910.562
$ws.Range("I62").Value = 11582.833
$ws.Range("K62").Value = 11582.833
$ws.Range("M62").Value = -10958.833
$ws.Range("H65").Value = 19910.562
$ws.Range("I65").Value = 11582.833
$ws.Range("K65").Value = 57914.165
$ws.Range("M65").Value = -54794.165
$ws.Range("H99").Value = 13433.689
$ws.Range("I99").Value = 7845.3887
$ws.Range("K99").Value = 7845.3887
$ws.Range("M99").Value = -6347.3887
$ws.Range("H126").Value = 13433.689
$ws.Range("I126").Value = 7845.3887
$ws.Range("K126").Value = 23536.1661
$ws.Range("M126").Value = -21066.1661
$ws.Range("H132").Value = 1428.963
$ws.Range("I132").Value = 1446.2307
$ws.Range("K132").Value = 4338.6921
$ws.Range("M132").Value = -1808.6921
$ws.Range("H136").Value = 2471
$ws.Range("I136").Value = 1760.8823
$ws.Range("K136").Value = 5282.6469
$ws.Range("M136").Value = -2732.6469

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 55555700
$ws.Range("J23").Value = 66666830
$ws.Range("L23").Value = 200000490
$ws.Range("N23").Value = -200000960
$ws.Range("H40").Value = 204.25
$ws.Range("J40").Value = 313.16666
$ws.Range("L40").Value = 1252.66664
$ws.Range("N40").Value = -1390.66664
$ws.Range("H70").Value = 11590.143
$ws.Range("I70").Value = 5559.8
$ws.Range("K70").Value = 16679.4
$ws.Range("M70").Value = -16364.4
$ws.Range("H73").Value = 11590.143
$ws.Range("I73").Value = 5559.8
$ws.Range("K73").Value = 16679.4
$ws.Range("M73").Value = -15587.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 99999
$ws.Range("J88").Value = 99999
$ws.Range("L88").Value = 99999
$ws.Range("N88").Value = -100901
$ws.Range("H91").Value = 99999
$ws.Range("J91").Value = 99999
$ws.Range("L91").Value = 99999
$ws.Range("N91").Value = -103119
$ws.Range("H102").Value = 2352.1
$ws.Range("I102").Value = 1890.5
$ws.Range("K102").Value = 1890.5
$ws.Range("M102").Value = -268.5
$ws.Range("H108").Value = 119999.664
$ws.Range("J108").Value = 119999.664
$ws.Range("L108").Value = 119999.664
$ws.Range("N108").Value = -127679.664
$ws.Range("H124").Value = 99999
$ws.Range("J124").Value = 99999
$ws.Range("L124").Value = 99999
$ws.Range("N124").Value = -109819
$ws.Range("H138").Value = 91731
$ws.Range("I138").Value = 50390
$ws.Range("J138").Value = 99999.2
$ws.Range("K138").Value = 50390
$ws.Range("L138").Value = 99999.2
$ws.Range("M138").Value = -45250
$ws.Range("N138").Value = -110279.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5221
$ws.Range("I40").Value = 4202.3335
$ws.Range("J40").Value = 9499.4
$ws.Range("K40").Value = 4202.3335
$ws.Range("L40").Value = 9499.4
$ws.Range("M40").Value = -4066.3335
$ws.Range("N40").Value = -9771.4
$ws.Range("H46").Value = 994.625
$ws.Range("I46").Value = 786
$ws.Range("J46").Value = 1119.8
$ws.Range("K46").Value = 786
$ws.Range("L46").Value = 1119.8
$ws.Range("M46").Value = -598
$ws.Range("N46").Value = -1495.8
$ws.Range("H55").Value = 1323.091
$ws.Range("I55").Value = 1013.2222
$ws.Range("K55").Value = 1013.2222
$ws.Range("M55").Value = -840.2222
$ws.Range("H68").Value = 3475317.5
$ws.Range("I68").Value = 6945994.5
$ws.Range("K68").Value = 6945994.5
$ws.Range("M68").Value = -6945245.5
$ws.Range("H71").Value = 3475317.5
$ws.Range("I71").Value = 6945994.5
$ws.Range("K71").Value = 34729972.5
$ws.Range("M71").Value = -34726228.5
$ws.Range("H122").Value = 3456.2
$ws.Range("I122").Value = 3341.491
$ws.Range("J122").Value = 4718
$ws.Range("K122").Value = 10024.473
$ws.Range("L122").Value = 14154
$ws.Range("M122").Value = -7574.473
$ws.Range("N122").Value = -19054
$ws.Range("H132").Value = 4864.654
$ws.Range("I132").Value = 2232.6428
$ws.Range("J132").Value = 7935.3335
$ws.Range("K132").Value = 6697.928400000001
$ws.Range("L132").Value = 23806.0005
$ws.Range("M132").Value = -4167.928400000001
$ws.Range("N132").Value = -28866.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 126666.664
$ws.Range("J4").Value = 126666.664
$ws.Range("L4").Value = 126666.664
$ws.Range("N4").Value = -126892.664
$ws.Range("H107").Value = 3631.1794
$ws.Range("I107").Value = 2120.5652
$ws.Range("K107").Value = 6361.6956
$ws.Range("M107").Value = -4441.6956
$ws.Range("H113").Value = 659.5
$ws.Range("I113").Value = 560
$ws.Range("K113").Value = 1680
$ws.Range("M113").Value = 490
$ws.Range("H122").Value = 1980.9062
$ws.Range("I122").Value = 1809.6666
$ws.Range("K122").Value = 5428.9998
$ws.Range("M122").Value = -2978.9998
$ws.Range("H136").Value = 253381.05
$ws.Range("I136").Value = 3650.303
$ws.Range("K136").Value = 10950.909
$ws.Range("M136").Value = -8400.909
